$d = $word.ActiveDocument

# The "Road buffer and setback requirements" illustration is a tiny inline
# picture (1x1 placeholder). The edit replaces it with a plain hyperlink
# whose visible text is the image's public URL, styled with the
# "Hyperlink" character style - i.e. the picture becomes a text link.

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F04_Road_Buffer_and_Setback.jpg?h=100%25&w=100%25"

$shp = $d.InlineShapes.Item(1)
$picRange = $shp.Range
$insertAt = $picRange.Start

# Remove the picture run entirely, then drop a hyperlink in its place.
$shp.Delete()

$target = $d.Range($insertAt, $insertAt)
$link = $d.Hyperlinks.Add($target, $url, "", "", $url)
